$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26/27 swap: Monero <-> LidoDAOToken (B and C columns) ---
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"

# --- Column D (Price) updates: force text to preserve exact formatting ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.391.84"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.824.31"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.50"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5328"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4039"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07622"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.108"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.325"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.000"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.600"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.85"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.828.92"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001076"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.44"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06597"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.079"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.397.21"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.17"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.205"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.58"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.452"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.62"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.039.19"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.98"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.123"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1102"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07413"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.646"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2232"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.901"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.30"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6253"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.176"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.395"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.53"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.701"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5837"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.99"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.990"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.199"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06889"

# --- Column E (Volume/1h) updates ---
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("E8").Value = "  +7.09%  "
$ws.Range("E9").Value = "  +2.46%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("E12").Value = "  +4.07%  "
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("E14").Value = "  +5.42%  "
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("E16").Value = "  +3.08%  "
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("E22").Value = "  +3.45%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("E25").Value = "  +5.59%  "
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("E27").Value = "  +7.10%  "
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("E29").Value = "  +3.00%  "
$ws.Range("E30").Value = "  +3.43%  "
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("E32").Value = "  +5.05%  "
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("E34").Value = "  +15.44%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("E38").Value = "  +5.59%  "
$ws.Range("E39").Value = "  +4.54%  "
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("E46").Value = "  +1.04%  "
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("E49").Value = "  +3.57%  "
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("E51").Value = "  +1.49%  "
